$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": move the saved selection from A7:XFD13 to B5.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Select()

# ---------------------------------------------------------------------------
# Sheet "Repayment schedule": remove column O's cells for rows 2-14
# (clear formatting + value so the <c> node is dropped entirely on save).
# ---------------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
for ($row = 2; $row -le 14; $row++) {
    $cell = $wsRepay.Cells.Item($row, 15)  # column O = 15
    $cell.ClearFormats()
    $cell.Value = ""
}

# ---------------------------------------------------------------------------
# Sheet "Transactions": update figures in rows 2-4 and move the selection
# from C12 to A2 (selected/activated last so it ends up the active tab,
# matching the original workbook's saved state).
# ---------------------------------------------------------------------------
$wsTxn = $wb.Worksheets.Item("Transactions")

$wsTxn.Range("A2").Value = 36
$wsTxn.Range("J2").Value = 9133.2199999999993

$wsTxn.Range("A3").Value = 34
$wsTxn.Range("C3").Value = 42064
$wsTxn.Range("E3").Value = 963.77
$wsTxn.Range("F3").Value = 866.78
$wsTxn.Range("G3").Value = 96.99
$wsTxn.Range("J3").Value = 4133.22

$wsTxn.Range("A4").Value = 32

$wsTxn.Range("A2").Select()
